# Applies the "Add files via upload" edit described by the commit diff:
#  - Channel sheet: update two IP addresses, adjust the view/selection
#  - Frame sheet: tweak a few start_address / read_byte numbers, move selection
#  - Details sheet: renumber several start_address / data_type values and
#    append three new rows (room2 sensor: gas / water / power usage)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Channel sheet (sheet1)
# ---------------------------------------------------------------------
$wsChannel = $wb.Worksheets.Item("Channel")
$wsChannel.Range("D2").Value = "127.0.0.1"
$wsChannel.Range("D3").Value = "192.168.0.58"
[void]$wsChannel.Range("E3").Select()

# ---------------------------------------------------------------------
# Details sheet (sheet3)
# ---------------------------------------------------------------------
$wsDetails = $wb.Worksheets.Item("Details")
$wsDetails.Range("J2").Value = 0

$wsDetails.Range("J3").Value = 2

$wsDetails.Range("H4").Value = 504
$wsDetails.Range("J4").Value = 3

$wsDetails.Range("H5").Value = 506
$wsDetails.Range("J5").Value = 4

$wsDetails.Range("H6").Value = 508
$wsDetails.Range("J6").Value = 5

$wsDetails.Range("H7").Value = 512
$wsDetails.Range("J7").Value = 0

$wsDetails.Range("H8").Value = 513

$wsDetails.Range("H9").Value = 514

$wsDetails.Range("H10").Value = 0

# New rows for room2 sensor1 usage channels
$wsDetails.Range("A11").Value = 10
$wsDetails.Range("B11").Value = "room2.센서1.가스량"
$wsDetails.Range("C11").Value = 2
$wsDetails.Range("D11").Value = 3
$wsDetails.Range("E11").Value = "AI"
$wsDetails.Range("F11").Value = 0
$wsDetails.Range("G11").Value = 2000
$wsDetails.Range("H11").Value = 1
$wsDetails.Range("I11").Value = 0
$wsDetails.Range("J11").Value = 1
$wsDetails.Range("K11").Value = 1
$wsDetails.Range("L11").Value = 1
$wsDetails.Range("M11").Value = 1

$wsDetails.Range("A12").Value = 11
$wsDetails.Range("B12").Value = "room2.센서1.수도량"
$wsDetails.Range("C12").Value = 2
$wsDetails.Range("D12").Value = 3
$wsDetails.Range("E12").Value = "AI"
$wsDetails.Range("F12").Value = 0
$wsDetails.Range("G12").Value = 2000
$wsDetails.Range("H12").Value = 2
$wsDetails.Range("I12").Value = 0
$wsDetails.Range("J12").Value = 1
$wsDetails.Range("K12").Value = 1
$wsDetails.Range("L12").Value = 1
$wsDetails.Range("M12").Value = 1

$wsDetails.Range("A13").Value = 12
$wsDetails.Range("B13").Value = "room2.센서1.전력량"
$wsDetails.Range("C13").Value = 2
$wsDetails.Range("D13").Value = 3
$wsDetails.Range("E13").Value = "AI"
$wsDetails.Range("F13").Value = 0
$wsDetails.Range("G13").Value = 2000
$wsDetails.Range("H13").Value = 3
$wsDetails.Range("I13").Value = 0
$wsDetails.Range("J13").Value = 1
$wsDetails.Range("K13").Value = 1
$wsDetails.Range("L13").Value = 1
$wsDetails.Range("M13").Value = 1

[void]$wsDetails.Range("A13").Select()

# ---------------------------------------------------------------------
# Frame sheet (sheet2) - select last so it stays the active tab
# ---------------------------------------------------------------------
$wsFrame = $wb.Worksheets.Item("Frame")
$wsFrame.Range("F2").Value = 501
$wsFrame.Range("G2").Value = 14
$wsFrame.Range("F3").Value = 0
$wsFrame.Range("F4").Value = 0
[void]$wsFrame.Range("F4").Select()
